$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$ws.Range("A5").Value = "177. Nth Highest Salary"

$ws.Range("B5").Value = "Medium"
$ws.Range("B5").Interior.Color = 49407

$ws.Range("C5").Value = "Data Manipulation"

$ws.Range("E5").Value = "https://leetcode.com/problems/nth-highest-salary/solutions/3858402/very-simple-and-clean-pandas-with-comments/comments/2187013 "
$ws.Hyperlinks.Add($ws.Range("E5"), "https://leetcode.com/problems/nth-highest-salary/solutions/3858402/very-simple-and-clean-pandas-with-comments/comments/2187013 ")
$ws.Range("E5").Style = $ws.Range("E4").Style

$ws.Range("D5").Value = ".drop_duplicates(), then sort_values(ascending=False), then if N exceeds num of unique salaries, or is 0 or less, return none, get the nth highest with .iloc[N-1]. Return the dynamic column with return pd.DataFrame({f'getNthHighestSalary({N})': [nth_highest]})."

$ws.Range("D16").Select()
